$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "(modified)" tag from column F for the existing cards
# (these cards are no longer considered "modified" relative to the original).
$clearRows = 2,3,4,5,7,11,13,15,16,17,19,20,21,22
foreach ($r in $clearRows) {
    $ws.Range("F" + $r).ClearContents()
}

# Add new Insanity cards from Horrific Journeys

$ws.Range("A23").Value = 'Paralyzing Fear'
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = 'Get close to -- that?! No! Absolutely not!'
$ws.Range("D23").Value = '<p>You cannot voluntarily move into a space that contains 1 or more monsters. You win or lose the game as normal.'
$ws.Range("F23").Value = 'same as original card'
$ws.Rows.Item(23).RowHeight = 30

$ws.Range("A24").Value = 'Right Time, Right Place'
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 'You just want to be home.'
$ws.Range("D24").Value = '<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you are not Lost in Time and Space. Otherwise, you lose the game.'
$ws.Range("F24").Value = 'same as original card'
$ws.Rows.Item(24).RowHeight = 45

$ws.Range("A25").Value = 'Deep One Hybrid'
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = 'Your ancestors trace back to Y''ha-nthlei and Innsmouth. They demand that you prove your worth.'
$ws.Range("D25").Value = '<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you have 2 or more <em>Evidence</em>.</p>' + "`n" + '<p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent''s, you take one possession.</p>'
$ws.Range("F25").Value = 'removed traitor aspect - come up with something better?'
$ws.Rows.Item(25).RowHeight = 75

# Update the frozen-pane view + selection to reflect the new content
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("D25").Select()
